$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "32.189272,34.881159"
$ws.Range("B28").Value = "32.706861,35.173861"
$ws.Range("B38").Value = "31.750585,35.215673"
$ws.Range("B64").Value = "32.018460,34.748167"
$ws.Range("B85").Value = "31.945204,34.878075"
$ws.Range("B91").Value = "31.755957,34.989832"
$ws.Range("B93").Value = "31.750898,35.207819"
$ws.Range("B95").Value = "31.858601,35.215336"
$ws.Range("B100").Value = "31.226237,34.809557"
$ws.Range("B101").Value = "31.068012,35.007848"
$ws.Range("B102").Value = "31.238084,34.794545"
$ws.Range("B105").Value = "31.863239,34.743120"
$ws.Range("B106").Value = "31.928344,34.878259"
$ws.Range("B107").Value = "32.174304,34.930966"
$ws.Range("B108").Value = "32.045852,34.752438"
$ws.Range("B113").Value = "31.663407,34.599960"
$ws.Range("B120").Value = "32.093309,34.885509"
$ws.Range("B121").Value = "31.246177,34.808709"
$ws.Range("B123").Value = "31.756796,34.988601"
